# Add season-record columns (Wins / Losses / Ties) to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new labels in AD1:AF1 ---
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the existing header formatting (bold, centered, thin border) by
# copying the format from the existing header cell A1 onto the new ones.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# --- Data rows (2-61): every player gets the team's season record ---
$wins = 66
$losses = 96
$ties = 0

$lastRow = 61
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = $wins
    $ws.Cells.Item($r, 31).Value = $losses
    $ws.Cells.Item($r, 32).Value = $ties
}
